# Daily attendance processing - reverse the order of names/emails listed
# in the "Recorded By" column (G) for every data row on the active sheet.
#
# Example: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# Cells that contain only a single value (no comma) are left untouched,
# since reversing a one-element list is a no-op.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row so we cover the whole "Recorded By" column.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ', '

        if ($parts.Length -gt 1) {
            $reversed = $parts[($parts.Length - 1)..0]
            $newVal = $reversed -join ', '
            $cell.Value = $newVal
        }
    }
}
